$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 159, shifting existing rows 159-162 down to 160-163
$ws.Rows.Item(159).Insert()

# Fill in the new row 159 with data
$ws.Range("A159").Value = 5
$ws.Range("B159").Value = "Macroferia Regional de Talca"
$ws.Range("C159").Value = "Maule"
$ws.Range("D159").Value = 44448
$ws.Range("D159").NumberFormat = $ws.Range("D160").NumberFormat
$ws.Range("E159").Value = 7
$ws.Range("F159").Value = 100112009
$ws.Range("G159").Value = "Acelga"
$ws.Range("H159").Value = "Sin especificar"
$ws.Range("I159").Value = "Primera"
$ws.Range("J159").Value = 400
$ws.Range("K159").Value = 2500
$ws.Range("L159").Value = 2500
$ws.Range("M159").Value = 2500
$ws.Range("N159").Value = "`$/docena de atados (4 kilos)"
$ws.Range("O159").Value = "Región del Maule"
$ws.Range("P159").Value = 625
$ws.Range("Q159").Value = 4
$ws.Range("R159").Value = "Hortaliza"
